$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to lowercase/underscore convention, consolidating
# "Fiscal Year + date" naming into a single consistent variable-name style.
$ws.Range("B1").Value = "actual"
$ws.Range("C1").Value = "actual_lastweek"
$ws.Range("D1").Value = "actual_lastyear"
$ws.Range("E1").Value = "target"
$ws.Range("A1").Value = "indicator_name"

# Update the selected range to match the header row selection.
$ws.Range("A1:E1").Select()
